$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.042.38"
$ws.Range("E2").Value = "  -1.87%  "

# Row 3
$ws.Range("D3").Value = "1.641.28"
$ws.Range("E3").Value = "  -1.89%  "

# Row 4
$ws.Range("E4").Value = "  +1.35%  "

# Row 5
$ws.Range("D5").Value = "'216.73"
$ws.Range("E5").Value = "  -1.37%  "

# Row 6
$ws.Range("D6").Value = "'0.5071"
$ws.Range("E6").Value = "  -1.81%  "

# Row 7
$ws.Range("D7").Value = "'1.021"
$ws.Range("E7").Value = "  +1.46%  "

# Row 8
$ws.Range("D8").Value = "'0.2588"
$ws.Range("E8").Value = "  +0.07%  "

# Row 9
$ws.Range("D9").Value = "'0.06431"
$ws.Range("E9").Value = "  -0.87%  "

# Row 10
$ws.Range("D10").Value = "'19.56"
$ws.Range("E10").Value = "  -2.72%  "

# Row 11
$ws.Range("D11").Value = "'0.07778"
$ws.Range("E11").Value = "  +1.35%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.660.80"
$ws.Range("E12").Value = "  -1.24%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.272"
$ws.Range("E13").Value = "  -1.85%  "

# Row 14
$ws.Range("D14").Value = "1.866.97"
$ws.Range("E14").Value = "  -1.89%  "

# Row 15
$ws.Range("D15").Value = "'0.5478"
$ws.Range("E15").Value = "  -2.21%  "

# Row 16
$ws.Range("D16").Value = "0.0₅7975"
$ws.Range("E16").Value = "  -1.05%  "

# Row 17
$ws.Range("D17").Value = "'63.68"
$ws.Range("E17").Value = "  -2.00%  "

# Row 18
$ws.Range("D18").Value = "26.055.67"
$ws.Range("E18").Value = "  -1.89%  "

# Row 19
$ws.Range("D19").Value = "'1.022"
$ws.Range("E19").Value = "  +1.48%  "

# Row 20
$ws.Range("D20").Value = "'204.96"
$ws.Range("E20").Value = "  -3.19%  "

# Row 21
$ws.Range("D21").Value = "'4.332"
$ws.Range("E21").Value = "  -2.74%  "

# Row 22
$ws.Range("D22").Value = "'10.04"
$ws.Range("E22").Value = "  -1.10%  "

# Row 23
$ws.Range("E23").Value = "  +0.90%  "

# Row 24
$ws.Range("E24").Value = "  +1.52%  "

# Row 25
$ws.Range("D25").Value = "'1.999"
$ws.Range("E25").Value = "  +15.47%  "

# Row 26
$ws.Range("D26").Value = "'142.55"
$ws.Range("E26").Value = "  -1.29%  "

# Row 27
$ws.Range("D27").Value = "'0.1160"
$ws.Range("E27").Value = "  -0.91%  "

# Row 28
$ws.Range("D28").Value = "'15.77"
$ws.Range("E28").Value = "  -0.07%  "

# Row 29
$ws.Range("D29").Value = "'6.831"
$ws.Range("E29").Value = "  -2.79%  "

# Row 30
$ws.Range("D30").Value = "'1.249"
$ws.Range("E30").Value = "  -1.26%  "

# Row 31
$ws.Range("D31").Value = "'0.05005"
$ws.Range("E31").Value = "  -4.32%  "

# Row 32
$ws.Range("D32").Value = "'3.274"
$ws.Range("E32").Value = "  -2.97%  "

# Row 33
$ws.Range("D33").Value = "'3.218"
$ws.Range("E33").Value = "  -0.29%  "

# Row 34
$ws.Range("D34").Value = "'1.542"
$ws.Range("E34").Value = "  -2.93%  "

# Row 35
$ws.Range("D35").Value = "'2.361"
$ws.Range("E35").Value = "  -0.53%  "

# Row 36
$ws.Range("D36").Value = "'2.657"
$ws.Range("E36").Value = "  -4.25%  "

# Row 37
$ws.Range("D37").Value = "'0.8959"
$ws.Range("E37").Value = "  -3.49%  "

# Row 38
$ws.Range("D38").Value = "'0.5679"

# Row 39
$ws.Range("D39").Value = "1.122.74"
$ws.Range("E39").Value = "  -3.24%  "

# Row 40
$ws.Range("D40").Value = "'0.01566"
$ws.Range("E40").Value = "  -2.20%  "

# Row 41
$ws.Range("B41").Value = "mCoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D41").Value = "'2.602"
$ws.Range("E41").Value = "  +0.82%  "

# Row 42
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "'1.022"
$ws.Range("E42").Value = "  +1.65%  "

# Row 43
$ws.Range("E43").Value = "  +0.17%  "

# Row 44
$ws.Range("D44").Value = "'0.8189"
$ws.Range("E44").Value = "  -5.78%  "

# Row 45
$ws.Range("D45").Value = "'99.99"
$ws.Range("E45").Value = "  -0.41%  "

# Row 46
$ws.Range("D46").Value = "1.774.99"
$ws.Range("E46").Value = "  -2.02%  "

# Row 47
$ws.Range("E47").Value = "  -1.49%  "

# Row 48
$ws.Range("D48").Value = "'0.4570"
$ws.Range("E48").Value = "  +1.72%  "

# Row 49
$ws.Range("D49").Value = "'1.020"
$ws.Range("E49").Value = "  +1.32%  "

# Row 50
$ws.Range("D50").Value = "'55.01"
$ws.Range("E50").Value = "  -2.07%  "

# Row 51
$ws.Range("D51").Value = "'0.05050"
$ws.Range("E51").Value = "  -1.75%  "

